# Commit: "Added Word Output Name Variable and Changed to Docx"
#
# The first paragraph currently reads:
#   Name: [bookmark "NameField"]Folsom High School[/bookmark]
# with both runs explicitly sized at 10pt (w:sz 20 half-points).
#
# Target state:
#   Name: McDonalds Folsom
# with the explicit font-size override removed and the bookmark removed
# (it is being replaced elsewhere by a real mail-merge / output "Name"
# variable, per the commit message), and the placeholder school name
# swapped for the new value.

$d = $word.ActiveDocument

$para = $d.Paragraphs(1)
$range = $d.Range($para.Range.Start, $para.Range.End - 1)

$newParagraphXml = '<w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma"/></w:rPr>' +
    '<w:t xml:space="preserve">Name: </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma"/></w:rPr>' +
    '<w:t>McDonalds Folsom</w:t></w:r>' +
    '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
    'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$range.InsertXML($packageXml)
